$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 571; everything below (old 571-587) shifts down to 573-589.
$ws.Rows("571:572").Insert()

# New row 571: Sandia, calidad "Primera", origen Peru, semana nueva (fecha serial 44509)
$ws.Cells.Item(571, 1).Value2  = 6
$ws.Cells.Item(571, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(571, 3).Value2  = "Metropolitana"
$ws.Cells.Item(571, 4).Value2  = 44509
$ws.Cells.Item(571, 5).Value2  = 13
$ws.Cells.Item(571, 6).Value2  = 100112028
$ws.Cells.Item(571, 7).Value2  = "Sandia"
$ws.Cells.Item(571, 8).Value2  = "Sin especificar"
$ws.Cells.Item(571, 9).Value2  = "Primera"
$ws.Cells.Item(571, 10).Value2 = 4900
$ws.Cells.Item(571, 11).Value2 = 850
$ws.Cells.Item(571, 12).Value2 = 900
$ws.Cells.Item(571, 13).Value2 = 878
$ws.Cells.Item(571, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(571, 15).Value2 = "Perú"
$ws.Cells.Item(571, 16).Value2 = 878
$ws.Cells.Item(571, 17).Value2 = 1
$ws.Cells.Item(571, 18).Value2 = "Hortaliza"

# New row 572: Sandia, calidad "Segunda", origen Peru, misma semana nueva
$ws.Cells.Item(572, 1).Value2  = 6
$ws.Cells.Item(572, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(572, 3).Value2  = "Metropolitana"
$ws.Cells.Item(572, 4).Value2  = 44509
$ws.Cells.Item(572, 5).Value2  = 13
$ws.Cells.Item(572, 6).Value2  = 100112028
$ws.Cells.Item(572, 7).Value2  = "Sandia"
$ws.Cells.Item(572, 8).Value2  = "Sin especificar"
$ws.Cells.Item(572, 9).Value2  = "Segunda"
$ws.Cells.Item(572, 10).Value2 = 1900
$ws.Cells.Item(572, 11).Value2 = 700
$ws.Cells.Item(572, 12).Value2 = 700
$ws.Cells.Item(572, 13).Value2 = 700
$ws.Cells.Item(572, 14).Value2 = "$/kilo (volumen en unidades)"
$ws.Cells.Item(572, 15).Value2 = "Perú"
$ws.Cells.Item(572, 16).Value2 = 700
$ws.Cells.Item(572, 17).Value2 = 1
$ws.Cells.Item(572, 18).Value2 = "Hortaliza"
